$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 17 (shifts rows 17:135 down to 18:136)
$ws.Rows.Item(17).Insert()

# Populate the newly inserted row 17 with the new weekly data point
$ws.Range("A17").Value = 6
$ws.Range("B17").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C17").Value = "Metropolitana"
$ws.Range("D17").Value = 44558
$ws.Range("D17").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E17").Value = 13
$ws.Range("F17").Value = 100112029
$ws.Range("G17").Value = "Orégano"
$ws.Range("H17").Value = "Sin especificar"
$ws.Range("I17").Value = "Primera"
$ws.Range("J17").Value = 36
$ws.Range("K17").Value = 9500
$ws.Range("L17").Value = 10000
$ws.Range("M17").Value = 9736
$ws.Range("N17").Value = "$/docena de atados"
$ws.Range("O17").Value = "Región Metropolitana"
$ws.Range("P17").Value = 3245
$ws.Range("Q17").Value = 3
$ws.Range("R17").Value = "Hortaliza"
